# Formed the consolidated report
# Recompute the "Absent" column (H) from the "Real" column (E):
# a row counts as Absent (1) when the student was not marked Real present (E=0),
# and not Absent (0) when the student was marked Real present (E=1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 3; $row -le 21; $row++) {
    $real = $ws.Cells.Item($row, 5).Value()
    if ($real -eq 1) {
        $ws.Cells.Item($row, 8).Value = 0
    } else {
        $ws.Cells.Item($row, 8).Value = 1
    }
}
